# add : switch show tree on tag item page
# Adds two new child tag rows ("images_aeriennes_1" / "images_aeriennes_2")
# under the existing "images_aeriennes" tag, growing the Tableau1 table
# from A1:E41 to A1:E43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table by two rows (this expands the table ref / autofilter
# and the sheet dimension automatically).
$newRow1 = $lo.ListRows.Add()
$newRow2 = $lo.ListRows.Add()

# Fill column A (id) for both new rows first ...
$ws.Range("A42").Value = "images_aeriennes_1"
$ws.Range("A43").Value = "images_aeriennes_2"

# ... then column C (name) for both new rows ...
$ws.Range("C42").Value = "Images aériennes sous partie 1"
$ws.Range("C43").Value = "Images aériennes sous partie 2"

# ... then column B (parent_id), pointing both rows at "images_aeriennes".
$ws.Range("B42").Value = "images_aeriennes"
$ws.Range("B43").Value = "images_aeriennes"

# Column D (description) is intentionally left blank for these rows.
# Touch column E (doc_ids) so the trailing empty cell is still materialised.
$ws.Range("E42").WrapText = $true
$ws.Range("E43").WrapText = $true

# Column B widened slightly to fit the new longer parent_id/id values.
$ws.Columns.Item(2).ColumnWidth = 13.5

# Restore the frozen-pane scroll position and active selection the way
# the workbook was left after the edit.
$ws.Activate()
$null = $ws.Range("C44").Select()

Write-Host "Added rows 42-43 to Tableau1; new range:" $lo.Range.Address()
